# Daily attendance processing - 2026-02-07 09:06:19 UTC
# Reorders the comma-separated names in the "Recorded By" column (G) for
# several rows of the attendance report worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Administrator, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad"
$ws.Range("G3").Value = "Administrator, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Majorelle Magdy"
$ws.Range("G4").Value = "Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy"
$ws.Range("G5").Value = "Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G6").Value = "Dr. Mohammad El-Tanany, Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Manar Montaser"
$ws.Range("G7").Value = "Dr. Lamiaa Ossama, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab"
$ws.Range("G8").Value = "Dr. Abeer Ragab, Dr. Nada Mohammad"
$ws.Range("G9").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G11").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G12").Value = "Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim"
$ws.Range("G13").Value = "Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Esraa Mostafa"
$ws.Range("G17").Value = "Dr. Esraa Samy, Dr. Mohammad Safwat"
$ws.Range("G19").Value = "Dr. Rania Ahmad Youssef, Dr. Mariam Toma Gerges"
$ws.Range("G27").Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Range("G30").Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Shorok Mohammad"
